$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the style of the previous date/time rows so we reuse existing style indices
$ws.Range("A11").Copy()
$ws.Range("A12").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("B11").Copy()
$ws.Range("B12").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("A12").Value = 44614
$ws.Range("B12").Value = 1.5
$ws.Range("C12").Value = "Updating positioning and documentation"

$ws.Range("C12").Select()
